$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy formatting from the
# existing header cell H1 (bold, centered, thin-bordered) so the new
# headers match the rest of the header row, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF)
$values = @(
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(7, 9),
    @(5, 6),
    @(3, 4),
    @(6, 7),
    @(2, 3),
    @(1, 2)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
